$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J: copy formatting from H1 (existing header style)
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 10).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data values for new columns
$iValues = @{
    2 = 8
    3 = 1
    4 = 1
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 1
    10 = 1
    11 = 6
}

$jValues = @{
    2 = 8
    3 = 4
    4 = 6
    5 = 5
    6 = 4
    7 = 4
    8 = 6
    9 = 3
    10 = 2
    11 = 6
}

foreach ($row in 2..11) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
